# Added sorting of keys to allow starring for estimated years
#
# The underlying data didn't change, but year labels that correspond to
# *estimated* (interpolated/extrapolated via formula) data points are now
# marked with a trailing "*" and stored as text labels instead of plain
# numbers, so that a later sort-by-key step can keep them distinguishable
# from the real observed years.

$wb = $excel.ActiveWorkbook

# --- Sheet "1_1": rows 6 (2007) & 7 (2008) are formula-derived -> star them
$ws = $wb.Worksheets.Item("1_1")
$ws.Range("A6").Value = "2007*"
$ws.Range("A7").Value = "2008*"
$ws.Activate()
$ws.Range("D52").Select()

# --- Sheet "1_2": rows 10 (2010) & 11 (2011) are formula-derived -> star them
$ws = $wb.Worksheets.Item("1_2")
$ws.Range("A10").Value = "2010*"
$ws.Range("A11").Value = "2011*"
$ws.Activate()
$ws.Range("G17").Select()

# --- Sheet "3_2": rows 10 (2010) & 11 (2011) are formula-derived -> star them
$ws = $wb.Worksheets.Item("3_2")
$ws.Range("A10").Value = "2010*"
$ws.Range("A11").Value = "2011*"
$ws.Activate()
$ws.Range("A12").Select()

# --- Sheet "4_2": rows 4 (2004), 9 (2009) & 10 (2010) are formula-derived -> star them
$ws = $wb.Worksheets.Item("4_2")
$ws.Range("A4").Value = "2004*"
$ws.Range("A9").Value = "2009*"
$ws.Range("A10").Value = "2010*"
$ws.Activate()
$ws.Range("F23").Select()

# --- Sheet "5_3": selection-only change
$ws = $wb.Worksheets.Item("5_3")
$ws.Activate()
$ws.Range("N42").Select()

# --- Sheet "8_1": all the interpolated/extrapolated years are starred
$ws = $wb.Worksheets.Item("8_1")
$ws.Range("A3").Value = "2004*"
$ws.Range("A4").Value = "2005*"
$ws.Range("A6").Value = "2007*"
$ws.Range("A7").Value = "2008*"
$ws.Range("A9").Value = "2010*"
$ws.Range("A10").Value = "2011*"
$ws.Range("A12").Value = "2013*"
$ws.Range("A13").Value = "2014*"
$ws.Range("A14").Value = "2015*"
$ws.Activate()
$ws.Range("G10").Select()

# --- "meta" sheet: selection-only change, last so it isn't the active tab
$ws = $wb.Worksheets.Item("meta")
$ws.Range("D32").Select()

# Restore "8_1" as the active / visible tab (matches activeTab=21 in the diff)
$ws = $wb.Worksheets.Item("8_1")
$ws.Activate()
$ws.Range("G10").Select()
